# This script rewrites the numeric/text result cells produced by the
# "Stationary generator" experiment so that the workbook reflects a
# re-run of the generator (per commit message: "volver a generar
# problemas cuadraticos y lineales").
#
# All of the changed cells in the original workbook store their values
# as plain text (shared strings) even though most of them look like
# numbers. Setting .Value directly with a numeric-looking string makes
# Excel auto-convert the cell to a true number, so we force a text
# number format before assigning, then clear the formatting again so
# the cell style stays at its original (default) index.

function Set-TextValue {
    param(
        $Range,
        [string]$Text
    )
    $Range.NumberFormat = "@"
    $Range.Value = $Text
    $Range.ClearFormats()
}

$wb = $excel.ActiveWorkbook

# --- Restricciones_del_lider ---
$ws = $wb.Worksheets.Item("Restricciones_del_lider")
Set-TextValue $ws.Range("A2") "0.049999999999998934 - x + y"
Set-TextValue $ws.Range("B2") "-0.049999999999998934"
Set-TextValue $ws.Range("D2") "0.4"

# --- Restricciones_del_follower ---
$ws = $wb.Worksheets.Item("Restricciones_del_follower")
Set-TextValue $ws.Range("A2") "34.04484536082474 + x - 3.367697594501718y"
Set-TextValue $ws.Range("B2") "-54.04484536082474"
Set-TextValue $ws.Range("D2") "0.55"
Set-TextValue $ws.Range("E2") "-9.9"
Set-TextValue $ws.Range("F2") "-9.8"
Set-TextValue $ws.Range("A3") "-27.504 + 1.9100000000000001y"
Set-TextValue $ws.Range("B3") "27.504"
Set-TextValue $ws.Range("D3") "0.45"
Set-TextValue $ws.Range("E3") "6.2"
Set-TextValue $ws.Range("F3") "1.6"
Set-TextValue $ws.Range("A4") "-40.0 + 1.1102230246251565e-16y"
Set-TextValue $ws.Range("D4") "0.55"
Set-TextValue $ws.Range("E4") "2.2"

# --- Punto_modificado ---
$ws = $wb.Worksheets.Item("Punto_modificado")
Set-TextValue $ws.Range("A2") "14.45"
Set-TextValue $ws.Range("B2") "14.4"

# --- Vector_bf ---
# NOTE: worksheet names "Vector_bf" and "Vector_BF" differ only by case, and
# Worksheets.Item(<name>) resolves case-insensitively (always landing on the
# first match). Use the 1-based sheet index instead to address each one.
$ws = $wb.Worksheets.Item(5)
if ($ws.Name -ne "Vector_bf") { throw "expected sheet 5 to be Vector_bf, got $($ws.Name)" }
Set-TextValue $ws.Range("A2") "-52.007266323024055"

# --- Vector_BF ---
$ws = $wb.Worksheets.Item(6)
if ($ws.Name -ne "Vector_BF") { throw "expected sheet 6 to be Vector_BF, got $($ws.Name)" }
Set-TextValue $ws.Range("A2") "-30.240000000000002"
Set-TextValue $ws.Range("A3") "-77.66220618556702"

# --- Vector_Alpha --- (this one stays a real number, not text)
$ws = $wb.Worksheets.Item("Vector_Alpha")
$ws.Range("A2").Value = 2.91
